$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'29.050.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.24%  '

$ws.Range("D3").Value = "'1.902.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.83%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = "'333.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.37%  '

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").Value = "'0.4643"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.01%  '

$ws.Range("D8").Value = "'0.4120"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.81%  '

$ws.Range("D9").Value = "'47.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.17%  '

$ws.Range("D10").Value = "'0.08000"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.36%  '

$ws.Range("D11").Value = "'1.005"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.52%  '

$ws.Range("E12").Value = '  -0.83%  '

$ws.Range("D13").Value = "'1.902.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.05%  '

$ws.Range("D14").Value = "'5.940"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.68%  '

$ws.Range("D15").Value = "'7.091"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.42%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = "'89.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.63%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = "'1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("D19").Value = "'0.06588"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.83%  '

$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Value = "'29.108.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.39%  '

$ws.Range("D23").Value = "'5.441"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.60%  '

$ws.Range("E24").Value = '  +2.05%  '

$ws.Range("D25").Value = "'2.226"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.89%  '

$ws.Range("D26").Value = "'2.131.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.16%  '

$ws.Range("D27").Value = "'157.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.01%  '

$ws.Range("D28").Value = "'19.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.35%  '

$ws.Range("D29").Value = "'2.122"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.17%  '

$ws.Range("D30").Value = "'5.438"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.93%  '

$ws.Range("D31").Value = "'118.25"
$ws.Range("D31").Style = "Normal"

$ws.Range("E32").Value = '  +1.26%  '

$ws.Range("D33").Value = "'0.09408"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("E34").Value = '  +4.01%  '

$ws.Range("E35").Value = '  +0.24%  '

$ws.Range("D36").Value = "'5.304"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.97%  '

$ws.Range("D37").Value = "'0.06101"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("D38").Value = "'0.02245"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.12%  '

$ws.Range("D39").Value = "'8.362"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.44%  '

$ws.Range("E40").Value = '  -0.50%  '

$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D42").Value = "'0.5799"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.29%  '

$ws.Range("D43").Value = "'10.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.43%  '

$ws.Range("D44").Value = "'0.1826"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.69%  '

$ws.Range("E45").Value = '  -1.25%  '

$ws.Range("D46").Value = "'2.322"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +12.76%  '

$ws.Range("E47").Value = '  -1.02%  '

$ws.Range("D48").Value = "'12.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.02%  '

$ws.Range("D49").Value = "'1.913"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.19%  '

$ws.Range("D50").Value = "'0.07049"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.82%  '

$ws.Range("D51").Value = "'46.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +18.33%  '
